{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst target = paragraphs.items[4];\ntarget.clear();\ntarget.style = \"Subtitle\";\ntarget.insertText(\"ASC1\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the \"Text Box 7\" shape (the ASC1 text box anchored as a floating\n# drawing in its own paragraph) and remove the drawing, replacing it with\n# plain paragraph text styled as \"Subtitle\".\n$shape = $null\nfor ($i = 1; $i -le $d.Shapes.Count; $i++) {\n    $candidate = $d.Shapes.Item($i)\n    if ($candidate.Name -eq \"Text Box 7\") {\n        $shape = $candidate\n        break\n    }\n}\n\n$para = $d.Paragraphs.Item(5)\n$range = $para.Range\n\nif ($shape -ne $null) {\n    $shape.Delete()\n}\n\n$range.Text = \"ASC1\"\n$para.Style = \"Subtitle\"\n"}
